$d = $word.ActiveDocument

$d.Content.Find.Execute("79÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "70÷4=", 2)
$d.Content.Find.Execute("69÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "71÷7=", 2)
$d.Content.Find.Execute("45÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "52÷4=", 2)
$d.Content.Find.Execute("93÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "97÷2=", 2)
$d.Content.Find.Execute("27÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "41÷2=", 2)
$d.Content.Find.Execute("66÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "74÷2=", 2)
$d.Content.Find.Execute("51÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "31÷2=", 2)
$d.Content.Find.Execute("45÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "42÷3=", 2)
$d.Content.Find.Execute("45÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "73÷9=", 2)
$d.Content.Find.Execute("25÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "71÷9=", 2)
$d.Content.Find.Execute("13÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "48÷8=", 2)
$d.Content.Find.Execute("23÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "56÷8=", 2)
$d.Content.Find.Execute("95÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "57÷4=", 2)
$d.Content.Find.Execute("49÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "24÷3=", 2)
$d.Content.Find.Execute("20÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "89÷7=", 2)
$d.Content.Find.Execute("10÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "80÷6=", 2)
$d.Content.Find.Execute("66÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "27÷2=", 2)
$d.Content.Find.Execute("52÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "69÷2=", 2)
$d.Content.Find.Execute("26÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "26÷9=", 2)
$d.Content.Find.Execute("44÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "75÷8=", 2)
$d.Content.Find.Execute("34÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "58÷5=", 2)
$d.Content.Find.Execute("71÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "10÷8=", 2)
$d.Content.Find.Execute("88÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "93÷3=", 2)
$d.Content.Find.Execute("62÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "34÷3=", 2)
$d.Content.Find.Execute("65÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "26÷5=", 2)
